$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A16").Value = "category[name]"
$ws.Range("B16").Value = 6
$ws.Range("C16").Value = 64
$ws.Range("F16").Value = "(?=.*?[a-zA-Z]+.*?)(.*)"

$ws.Range("D16").Select() | Out-Null
